# This workbook is an old file being re-saved/synced ("Syncing a bunch of old
# files"). Most of the underlying XML churn in the source diff (namespace
# bumps, extLst blocks, theme font metadata, per-row x14ac:dyDescent hints,
# calcPr/fileVersion bumps, etc.) is Excel-version save noise that isn't
# driven by any user action in the UI, so it isn't reproduced here. The
# concrete, user-visible edits captured by this script are:
#   1) Column width adjustments on columns B, C, E, F and G.
#   2) The final selection left on the sheet (H32:H36, active cell H32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B -> stored width 12.88671875 characters (bestFit)
$ws.Columns.Item(2).ColumnWidth = 12
# Column C -> stored width 103.5546875 characters (bestFit)
$ws.Columns.Item(3).ColumnWidth = 102.66666666666667
# Column E -> stored width 6.33203125 characters
$ws.Columns.Item(5).ColumnWidth = 5.5
# Column F -> stored width 9.44140625 characters
$ws.Columns.Item(6).ColumnWidth = 8.666666666666666
# Column G -> stored width 64.44140625 characters
$ws.Columns.Item(7).ColumnWidth = 63.666666666666664

# Leave the selection where the author left it when they last saved the file.
$ws.Range("H32:H36").Select()
